# Actualización automática 2025-06-01 08:00:06
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item(1)

# Column I (9) width: 11 -> 9 (OOXML width units; ColumnWidth = OOXML width - 5/6)
$ws1.Columns.Item(9).ColumnWidth = 9 - 5/6

# All monthly figures reset to 0 for every advisor row (rows 2-21, columns C..N)
$ws1.Range("C2:N21").Value = 0

# Row 22 summary labels ("X de 20") all become "0 de 20"
$ws1.Range("C22:N22").Value = "0 de 20"

# --- Sheet 2: "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item(2)

# Column widths shift: D 14->13, E 13->14, F 14->11 (OOXML width units)
$ws2.Columns.Item(4).ColumnWidth = 13 - 5/6
$ws2.Columns.Item(5).ColumnWidth = 14 - 5/6
$ws2.Columns.Item(6).ColumnWidth = 11 - 5/6

# Header months roll forward by one: febrero/marzo/abril/mayo -> marzo/abril/mayo/junio
$ws2.Range("C1").Value = "marzo"
$ws2.Range("D1").Value = "abril"
$ws2.Range("E1").Value = "mayo"
$ws2.Range("F1").Value = "junio"

# Data rolls left by one month column; the newest month (F) starts at 0
for ($r = 2; $r -le 22; $r++) {
    $oldD = $ws2.Cells.Item($r, 4).Value2
    $oldE = $ws2.Cells.Item($r, 5).Value2
    $oldF = $ws2.Cells.Item($r, 6).Value2

    $ws2.Cells.Item($r, 3).Value = $oldD
    $ws2.Cells.Item($r, 4).Value = $oldE
    $ws2.Cells.Item($r, 5).Value = $oldF
    $ws2.Cells.Item($r, 6).Value = 0
}
